$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use a scratch row far below the real data so whole-row Copy() round trips
# preserve the exact underlying cell types (text stays text, e.g. the
# "HH:MM" / "YYYY-MM-DD" strings in columns Y/Z/AA/AB don't get reinterpreted
# as real dates the way a Value2 array write would).
#
# NOTE: this engine's Range.Copy(destination) only overwrites cells where the
# source has content - cells that are blank in the source are left untouched
# in the destination rather than being blanked out (unlike real Excel). So
# every destination range is explicitly Clear()-ed immediately before the
# Copy() that is meant to fully repopulate it.
$scratchRow = 1000

function FullRow($sheet, $r) {
    return $sheet.Range("A" + $r + ":AY" + $r)
}

function Swap-Rows($sheet, $r1, $r2, $scratchR) {
    $row1 = FullRow $sheet $r1
    $row2 = FullRow $sheet $r2
    $scratch = FullRow $sheet $scratchR

    $scratch.Clear()
    $row1.Copy($scratch)

    $row1.Clear()
    $row2.Copy($row1)

    $row2.Clear()
    $scratch.Copy($row2)

    $scratch.Clear()
}

# Row 3 <-> Row 4 (full row content swap)
Swap-Rows $ws 3 4 $scratchRow

# Row 6 <-> Row 7 (full row content swap)
Swap-Rows $ws 6 7 $scratchRow

# Row 11 <-> Row 12 (full row content swap)
Swap-Rows $ws 11 12 $scratchRow

# Rows 14 -> 15 -> 16 -> 14 (3-way rotation):
#   new row14 = old row15, new row15 = old row16, new row16 = old row14
$row14 = FullRow $ws 14
$row15 = FullRow $ws 15
$row16 = FullRow $ws 16
$scratch = FullRow $ws $scratchRow

$scratch.Clear()
$row14.Copy($scratch)

$row14.Clear()
$row15.Copy($row14)

$row15.Clear()
$row16.Copy($row15)

$row16.Clear()
$scratch.Copy($row16)

$scratch.Clear()
